$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new blank column before column N (14) to make room for an
# additional "Variable Instalments" field - everything from N onward
# (Late/Outstanding/heading/Outstanding) shifts right by one column.
$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Mirror where the user ended up clicking after inserting the column.
$ws.Range("R6").Select()
